# Update "F" column (想去人数 / interested-count) values on each sheet to
# reflect a fresh data scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 76
$ws1.Range("F6").Value  = 2469
$ws1.Range("F7").Value  = 51
$ws1.Range("F9").Value  = 523
$ws1.Range("F10").Value = 1510
$ws1.Range("F11").Value = 19
$ws1.Range("F12").Value = 608
$ws1.Range("F13").Value = 1349
$ws1.Range("F14").Value = 1349
$ws1.Range("F15").Value = 1195
$ws1.Range("F16").Value = 493
$ws1.Range("F17").Value = 3509
$ws1.Range("F18").Value = 639
$ws1.Range("F19").Value = 3263
$ws1.Range("F21").Value = 609
$ws1.Range("F22").Value = 23
$ws1.Range("F23").Value = 281
$ws1.Range("F25").Value = 1108
$ws1.Range("F28").Value = 946
$ws1.Range("F29").Value = 930

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value  = 84
$ws2.Range("F10").Value = 15
$ws2.Range("F11").Value = 85
$ws2.Range("F18").Value = 109
$ws2.Range("F19").Value = 232
$ws2.Range("F21").Value = 465

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 487

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 76
$ws4.Range("F11").Value = 487
$ws4.Range("F12").Value = 2469
$ws4.Range("F14").Value = 51
$ws4.Range("F17").Value = 84
$ws4.Range("F19").Value = 523
$ws4.Range("F20").Value = 15
$ws4.Range("F21").Value = 85
$ws4.Range("F22").Value = 1510
$ws4.Range("F24").Value = 19
$ws4.Range("F25").Value = 1350
$ws4.Range("F26").Value = 1350
$ws4.Range("F29").Value = 1195
$ws4.Range("F30").Value = 493
$ws4.Range("F32").Value = 3509
$ws4.Range("F33").Value = 639
$ws4.Range("F34").Value = 3263
$ws4.Range("F37").Value = 609
$ws4.Range("F38").Value = 23
$ws4.Range("F39").Value = 281
$ws4.Range("F40").Value = 1108
$ws4.Range("F42").Value = 109
$ws4.Range("F43").Value = 232
$ws4.Range("F45").Value = 465
$ws4.Range("F48").Value = 946
$ws4.Range("F49").Value = 930
